# blacklisted_words.xlsx — "updated master data with Spanish lang"
#
# Before: rows 2-7 held
#   eng/shit, eng/damn, eng/nigga, eng/dammit, spa/Merde, spa/bon sang
# After: rows 2-3 hold just
#   eng/shit, spa/bon sang
# i.e. the English profanity variants (damn/nigga/dammit) and the old
# French "Merde" Spanish-row entry are dropped, replaced by a single,
# cleaned-up Spanish row; the trailing blank D-only rows shift up to
# follow directly after the new row 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 3:7 (damn / nigga / dammit / Merde-spa / bon sang-spa)
$ws.Rows("3:7").Delete()

# Insert a fresh row 3 to hold the updated Spanish blacklist entry
$ws.Rows("3:3").Insert()

$ws.Range("A3").Value = "spa"
$ws.Range("B3").Value = "bon sang"
$ws.Range("C3").Value = "Mot sur la liste noire"

# Copy D2 ("TRUE", stored as text via column D's text format) into D3 so it
# keeps the same shared-string/text representation instead of becoming a
# native boolean.
$ws.Range("D2").Copy($ws.Range("D3"))

# Match the saved selection state
[void]$ws.Range("B5").Select()
